$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like numbers; force text format so Excel
# does not silently convert them (e.g. "315.32" -> 315.32 numeric).
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D14", "D15", "D17", "D21", "D22", "D23", "D26", "D27", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D43", "D44", "D45", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (prices in column D, volume change % in column E,
# and the swapped Maker / ApeXProtocol rows 42-43 in columns B/C).
$ws.Range("D2").Value = "40.812.96"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.413.00"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "315.32"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "88.77"
$ws.Range("E6").Value = "  -3.93%  "
$ws.Range("D7").Value = "0.537"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -4.59%  "
$ws.Range("D10").Value = "31.83"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").Value = "0.0827"
$ws.Range("E11").Value = "  -4.62%  "
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("D13").Value = "2.795.63"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "2.416.74"
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").Value = "0.766"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "40.794.60"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("E19").Value = "  -3.68%  "
$ws.Range("E20").Value = "  -4.28%  "
$ws.Range("D21").Value = "71.38"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "10.99"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").Value = "234.60"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "23.95"
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("D30").Value = "34.79"
$ws.Range("E30").Value = "  -4.97%  "
$ws.Range("D31").Value = "155.41"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("D34").Value = "2.51"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").Value = "0.0741"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").Value = "2.96"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "16.61"
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").Value = "0.0998"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.986.22"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -8.47%  "
$ws.Range("D44").Value = "18.66"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("D45").Value = "0.0273"
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("E46").Value = "  -4.90%  "
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "2.657.84"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "94.36"
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("D50").Value = "72.78"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "51.62"

# Restore default cell style on the forced-text cells so we do not leave
# a lingering custom number format applied to them.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
